$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force Text format so values like "38.754.02" or
#     "0.612" are stored as literal text, matching the original inlineStr cells
#     instead of being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.754.02'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.085.66'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.612'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.18'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0844'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.395.49'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.99'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.85'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.797'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.077.22'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.689.69'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.46'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.10'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.84'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.53'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.46'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.70'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0608'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.90'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.541.42'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.63'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0923'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.69'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.11'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.283.96'

# --- Volume(1h) column (E): percentage-change strings, already safe as text
#     (they contain a leading/trailing double-space and a percent sign so
#     Excel keeps them as text automatically).
$ws.Range("E2").Value = '  +2.75%  '
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("E7").Value = '  +1.50%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("E10").Value = '  +1.04%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("E13").Value = '  +4.26%  '
$ws.Range("E14").Value = '  +2.38%  '
$ws.Range("E15").Value = '  +4.70%  '
$ws.Range("E16").Value = '  +0.52%  '
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("E18").Value = '  +2.71%  '
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("E21").Value = '  +1.24%  '
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("E27").Value = '  +2.43%  '
$ws.Range("E28").Value = '  +7.01%  '
$ws.Range("E29").Value = '  +13.09%  '
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("E32").Value = '  +5.21%  '
$ws.Range("E33").Value = '  +2.91%  '
$ws.Range("E34").Value = '  +3.50%  '
$ws.Range("E35").Value = '  +0.91%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  +1.52%  '
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("E39").Value = '  -0.06%  '
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("E41").Value = '  +5.47%  '
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("E43").Value = '  +2.73%  '
$ws.Range("E44").Value = '  -0.87%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("E46").Value = '  +8.92%  '
$ws.Range("E47").Value = '  +1.37%  '
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("E49").Value = '  +2.86%  '
$ws.Range("E50").Value = '  +0.56%  '
